# Update cryptos list with latest prices and 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.230.66'
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").Value = '1.894.54'
$ws.Range("E3").Value = '  -1.17%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.03'
$ws.Range("E5").Value = '  -2.64%  '

$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5064'
$ws.Range("E7").Value = '  -3.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4041'
$ws.Range("E8").Value = '  -0.99%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08293'
$ws.Range("E9").Value = '  -2.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.113'
$ws.Range("E10").Value = '  -0.94%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.25'
$ws.Range("E11").Value = '  -1.44%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.22'
$ws.Range("E12").Value = '  +7.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.404'
$ws.Range("E13").Value = '  -0.77%  '

$ws.Range("D14").Value = '1.880.83'
$ws.Range("E14").Value = '  -1.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.333'
$ws.Range("E15").Value = '  -0.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  +0.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.93'
$ws.Range("E17").Value = '  -2.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001101'
$ws.Range("E18").Value = '  -1.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06468'
$ws.Range("E19").Value = '  -3.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.44'
$ws.Range("E20").Value = '  +0.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.924'
$ws.Range("E22").Value = '  -1.40%  '

$ws.Range("D23").Value = '30.221.39'
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.30'
$ws.Range("E24").Value = '  -0.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.185'
$ws.Range("E25").Value = '  -1.35%  '

$ws.Range("D26").Value = '2.111.66'
$ws.Range("E26").Value = '  -1.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.63'
$ws.Range("E27").Value = '  +2.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.92'
$ws.Range("E28").Value = '  +0.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.272'
$ws.Range("E29").Value = '  -5.89%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.09'
$ws.Range("E30").Value = '  +0.13%  '

$ws.Range("E31").Value = '  +2.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1043'
$ws.Range("E32").Value = '  -2.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.009'
$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.705'
$ws.Range("E34").Value = '  +1.82%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02452'
$ws.Range("E35").Value = '  -1.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.318'
$ws.Range("E36").Value = '  +2.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06449'
$ws.Range("E37").Value = '  -2.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2158'
$ws.Range("E38").Value = '  -2.17%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.188'
$ws.Range("E39").Value = '  -3.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.604'
$ws.Range("E40").Value = '  -3.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6404'
$ws.Range("E41").Value = '  -2.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.39'
$ws.Range("E42").Value = '  -2.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.213'

$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.25'
$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5976'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.141'
$ws.Range("E47").Value = '  +2.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.641'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.70'
$ws.Range("E49").Value = '  -0.71%  '

$ws.Range("E50").Value = '  -2.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.81'
$ws.Range("E51").Value = '  -1.09%  '
